# Gal-Galr2.xlsx update: refresh TPM-derived NATMI edge stats, drop the ECs-sending rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-10 (sending cluster = ECs) are removed entirely.
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Galr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.48067
$ws.Range("H2").Value = 1.44201
$ws.Range("I2").Value = 0.1949338371837906
$ws.Range("J2").Value = 0.1949338371837907
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.165415
$ws.Range("N2").Value = 0.496245
$ws.Range("O2").Value = 0.06044768156291203
$ws.Range("P2").Value = 0.06044768156291203
$ws.Range("Q2").Value = 0.07951002805
$ws.Range("R2").Value = 0.71559025245
$ws.Range("S2").Value = 0.01178329851592232
$ws.Range("T2").Value = 0.01178329851592232

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Galr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.48067
$ws.Range("H3").Value = 1.44201
$ws.Range("I3").Value = 0.1949338371837906
$ws.Range("J3").Value = 0.1949338371837907
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.25579
$ws.Range("N3").Value = 6.76737
$ws.Range("O3").Value = 0.8243344049378915
$ws.Range("P3").Value = 0.8243344049378915
$ws.Range("Q3").Value = 1.0842905793
$ws.Range("R3").Value = 9.758615213699999
$ws.Range("S3").Value = 0.1606906686771599
$ws.Range("T3").Value = 0.1606906686771599

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Galr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.48067
$ws.Range("H4").Value = 1.44201
$ws.Range("I4").Value = 0.1949338371837906
$ws.Range("J4").Value = 0.1949338371837907
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3152936666666666
$ws.Range("N4").Value = 0.945881
$ws.Range("O4").Value = 0.1152179134991965
$ws.Range("P4").Value = 0.1152179134991965
$ws.Range("Q4").Value = 0.1515522067566666
$ws.Range("R4").Value = 1.36396986081
$ws.Range("S4").Value = 0.02245986999070845
$ws.Range("T4").Value = 0.02245986999070846

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Galr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.985141
$ws.Range("H5").Value = 5.955423
$ws.Range("I5").Value = 0.8050661628162092
$ws.Range("J5").Value = 0.8050661628162094
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.165415
$ws.Range("N5").Value = 0.496245
$ws.Range("O5").Value = 0.06044768156291203
$ws.Range("P5").Value = 0.06044768156291203
$ws.Range("Q5").Value = 0.328372098515
$ws.Range("R5").Value = 2.955348886635
$ws.Range("S5").Value = 0.04866438304698971
$ws.Range("T5").Value = 0.04866438304698971

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gal"
$ws.Range("C6").Value = "Galr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.985141
$ws.Range("H6").Value = 5.955423
$ws.Range("I6").Value = 0.8050661628162092
$ws.Range("J6").Value = 0.8050661628162094
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.25579
$ws.Range("N6").Value = 6.76737
$ws.Range("O6").Value = 0.8243344049378915
$ws.Range("P6").Value = 0.8243344049378915
$ws.Range("Q6").Value = 4.478061216389999
$ws.Range("R6").Value = 40.30255094750999
$ws.Range("S6").Value = 0.6636437362607315
$ws.Range("T6").Value = 0.6636437362607316

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gal"
$ws.Range("C7").Value = "Galr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.985141
$ws.Range("H7").Value = 5.955423
$ws.Range("I7").Value = 0.8050661628162092
$ws.Range("J7").Value = 0.8050661628162094
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3152936666666666
$ws.Range("N7").Value = 0.945881
$ws.Range("O7").Value = 0.1152179134991965
$ws.Range("P7").Value = 0.1152179134991965
$ws.Range("Q7").Value = 0.6259023847403332
$ws.Range("R7").Value = 5.633121462663
$ws.Range("S7").Value = 0.09275804350848807
$ws.Range("T7").Value = 0.09275804350848808

